$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered) from C1 to D1:E1
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Column D: EMPLEADO (tipo de empleado: PERMANENTE / TEMPORAL) -- filled first, top to bottom
$ws.Range("D1").Value = "EMPLEADO"
$tipoEmpleado = @("PERMANENTE", "TEMPORAL", "PERMANENTE", "PERMANENTE", "TEMPORAL", "PERMANENTE", "PERMANENTE", "TEMPORAL", "PERMANENTE", "PERMANENTE")
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $tipoEmpleado[$i]
}

# Column E: ESTADO CIVIL (CASADO / SOLTERO / DIVORCIADO) -- filled second, top to bottom
$ws.Range("E1").Value = "ESTADO CIVIL"
$estadoCivil = @("CASADO", "CASADO", "SOLTERO", "CASADO", "SOLTERO", "DIVORCIADO", "CASADO", "CASADO", "CASADO", "CASADO")
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $estadoCivil[$i]
}

# Column widths (target raw widths 14.21875 / 13.44140625; engine quantizes
# ColumnWidth to 1/6-character steps, so pick inputs that land on the nearest
# reachable values: 14.16666... and 13.5)
$ws.Columns.Item(4).ColumnWidth = 13.33
$ws.Columns.Item(5).ColumnWidth = 12.65

# Update selection to match final state
$ws.Range("E12").Select() | Out-Null

$wb.Save()
